$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "GROUP"
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 3).Value = "Grp1"
}

$ws.Columns.Item(3).ColumnWidth = 17.15

$ws.Range("C7").Select() | Out-Null
